# fix bug lich su
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: B6 30 -> 210
$ws.Range("B6").Value = 210

# Row 9: B9 150 -> 50, C9 3 -> 9, D9 1 -> 3, E9 2 -> 6
$ws.Range("B9").Value = 50
$ws.Range("C9").Value = 9
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = 6

# Row 9 history strings (F9, G9, H9) get extra entries appended
$ws.Range("F9").Value = ";0;0;1;1;0;1;0;0;0"
$ws.Range("G9").Value = ";14;33;13;13;10;13;31;0;4"
$ws.Range("H9").Value = ";-120.0;-120;+50.0;+100;-100;+1499700;-1499700;-1499700;-100"

# Update the active selection to F15
$ws.Range("F15").Select()
